$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from G1 (header style) onto H1, then set its value.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("H1").Value = "Save"

$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
